$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 3.037522333333333
$ws.Range("H2").Value = 9.112567
$ws.Range("I2").Value = 0.1153015356242242
$ws.Range("J2").Value = 0.1153015356242242
$ws.Range("M2").Value = 0.029424
$ws.Range("N2").Value = 0.08827199999999999
$ws.Range("O2").Value = 0.1473063425232919
$ws.Range("P2").Value = 0.1473063425232919
$ws.Range("Q2").Value = 0.08937605713599998
$ws.Range("R2").Value = 0.8043845142239999
$ws.Range("S2").Value = 0.01698464750012352
$ws.Range("T2").Value = 0.01698464750012352
$ws.Range("G3").Value = 3.037522333333333
$ws.Range("H3").Value = 9.112567
$ws.Range("I3").Value = 0.1153015356242242
$ws.Range("J3").Value = 0.1153015356242242
$ws.Range("O3").Value = 0.852693657476708
$ws.Range("P3").Value = 0.852693657476708
$ws.Range("Q3").Value = 0.5173599163803334
$ws.Range("R3").Value = 4.656239247423001
$ws.Range("S3").Value = 0.09831688812410068
$ws.Range("T3").Value = 0.09831688812410068
$ws.Range("I4").Value = 0.325995654495798
$ws.Range("J4").Value = 0.325995654495798
$ws.Range("M4").Value = 0.029424
$ws.Range("N4").Value = 0.08827199999999999
$ws.Range("O4").Value = 0.1473063425232919
$ws.Range("P4").Value = 0.1473063425232919
$ws.Range("Q4").Value = 0.2526957345759999
$ws.Range("R4").Value = 2.274261611183999
$ws.Range("S4").Value = 0.04802122754226276
$ws.Range("T4").Value = 0.04802122754226276
$ws.Range("I5").Value = 0.325995654495798
$ws.Range("J5").Value = 0.325995654495798
$ws.Range("O5").Value = 0.852693657476708
$ws.Range("P5").Value = 0.852693657476708
$ws.Range("S5").Value = 0.2779744269535352
$ws.Range("T5").Value = 0.2779744269535352
$ws.Range("I6").Value = 0.5587028098799778
$ws.Range("J6").Value = 0.5587028098799777
$ws.Range("M6").Value = 0.029424
$ws.Range("N6").Value = 0.08827199999999999
$ws.Range("O6").Value = 0.1473063425232919
$ws.Range("P6").Value = 0.1473063425232919
$ws.Range("Q6").Value = 0.433078830976
$ws.Range("R6").Value = 3.897709478784
$ws.Range("S6").Value = 0.08230046748090566
$ws.Range("T6").Value = 0.08230046748090565
$ws.Range("I7").Value = 0.5587028098799778
$ws.Range("J7").Value = 0.5587028098799777
$ws.Range("O7").Value = 0.852693657476708
$ws.Range("P7").Value = 0.852693657476708
$ws.Range("S7").Value = 0.4764023423990721
$ws.Range("T7").Value = 0.476402342399072
